# Update Clcf1-Lifr NATMI TPM values with the newly recomputed per-cluster
# ligand/receptor expression figures, and re-derive every dependent column
# (detection rate, specificity, and edge weight/specificity products).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New per-cluster figures (the "new tpm" values from the commit message).
# Total cell counts per cluster are unchanged; only expression values moved.
# ---------------------------------------------------------------------------
$totalCells = @{
    "ECs"               = 3
    "FAPs"              = 3
    "Inflammatory-Mac"  = 3
    "MuSCs"             = 2
    "Resolving-Mac"     = 3
}

$newExpressingCells = @{
    "ECs"               = 2
    "FAPs"              = 3
    "Inflammatory-Mac"  = 3
    "MuSCs"             = 2
    "Resolving-Mac"     = 3
}

$newLigandAvg = @{
    "ECs"               = 1.352356
    "FAPs"              = 3.197813
    "Inflammatory-Mac"  = 2.438989
    "MuSCs"             = 8.165625500000001
    "Resolving-Mac"     = 0.819627
}

$newLigandTotal = @{
    "ECs"               = 4.057068
    "FAPs"              = 9.593439
    "Inflammatory-Mac"  = 7.316967
    "MuSCs"             = 16.331251
    "Resolving-Mac"     = 2.458881
}

$newReceptorAvg = @{
    "ECs"               = 23.80409633333333
    "FAPs"              = 35.705903
    "Inflammatory-Mac"  = 33.79564933333334
    "MuSCs"             = 9.2924895
    "Resolving-Mac"     = 14.353493
}

$newReceptorTotal = @{
    "ECs"               = 71.41228899999999
    "FAPs"              = 107.117709
    "Inflammatory-Mac"  = 101.386948
    "MuSCs"             = 18.584979
    "Resolving-Mac"     = 43.060479
}

# Sums across all clusters, used for the "derived specificity" columns.
$sumLigandAvg = 0.0
$sumLigandTotal = 0.0
$sumReceptorAvg = 0.0
$sumReceptorTotal = 0.0
foreach ($k in $newLigandAvg.Keys) { $sumLigandAvg += $newLigandAvg[$k] }
foreach ($k in $newLigandTotal.Keys) { $sumLigandTotal += $newLigandTotal[$k] }
foreach ($k in $newReceptorAvg.Keys) { $sumReceptorAvg += $newReceptorAvg[$k] }
foreach ($k in $newReceptorTotal.Keys) { $sumReceptorTotal += $newReceptorTotal[$k] }

# ---------------------------------------------------------------------------
# Columns: A Sending cluster, D Target cluster,
#          E Ligand-expressing cells, F Ligand detection rate,
#          G Ligand avg expr, H Ligand total expr,
#          I Ligand specificity(avg), J Ligand specificity(total),
#          M Receptor avg expr, N Receptor total expr,
#          O Receptor specificity(avg), P Receptor specificity(total),
#          Q Edge avg weight, R Edge total weight,
#          S Edge avg specificity, T Edge total specificity
# ---------------------------------------------------------------------------
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $sending = $ws.Cells.Item($r, 1).Value()
    $target  = $ws.Cells.Item($r, 4).Value()

    $e = $newExpressingCells[$sending]
    $f = $e / $totalCells[$sending]

    $g = $newLigandAvg[$sending]
    $h = $newLigandTotal[$sending]
    $i = $g / $sumLigandAvg
    $j = $h / $sumLigandTotal

    $m = $newReceptorAvg[$target]
    $n = $newReceptorTotal[$target]
    $o = $m / $sumReceptorAvg
    $p = $n / $sumReceptorTotal

    $q = $g * $m
    $rr = $h * $n
    $s = $i * $o
    $t = $j * $p

    $ws.Cells.Item($r, 5).Value  = $e    # E
    $ws.Cells.Item($r, 6).Value  = $f    # F
    $ws.Cells.Item($r, 7).Value  = $g    # G
    $ws.Cells.Item($r, 8).Value  = $h    # H
    $ws.Cells.Item($r, 9).Value  = $i    # I
    $ws.Cells.Item($r, 10).Value = $j    # J
    $ws.Cells.Item($r, 13).Value = $m    # M
    $ws.Cells.Item($r, 14).Value = $n    # N
    $ws.Cells.Item($r, 15).Value = $o    # O
    $ws.Cells.Item($r, 16).Value = $p    # P
    $ws.Cells.Item($r, 17).Value = $q    # Q
    $ws.Cells.Item($r, 18).Value = $rr   # R
    $ws.Cells.Item($r, 19).Value = $s    # S
    $ws.Cells.Item($r, 20).Value = $t    # T
}
